$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: mark A31 as "MEM!" with red fill
$ws.Range("A31").Value = "MEM!"
$ws.Range("A31").Interior.Color = 255

# Row 32: A32 = 7252
$ws.Range("A32").Value = 7252

# Row 33: A33 = 4880
$ws.Range("A33").Value = 4880

# Row 34: mark A34 as "MEM!" with red fill
$ws.Range("A34").Value = "MEM!"
$ws.Range("A34").Interior.Color = 255

# Row 35: A35 = 3476
$ws.Range("A35").Value = 3476

# Update selection to match final cursor position
$ws.Range("J36").Select() | Out-Null
